# Applies the edits described by the commit:
#  - Bus connections: line resistances/lengths (E2, E3) 7500 -> 10000 (dependent C2/C3 formulas recalc)
#  - Generator data: capacity (E4) 15000 -> 20000 (dependent F4/G4/H4 formulas recalc)
#  - Selection/active-tab bookkeeping: "Bus connections" becomes the active sheet/tab,
#    "Bus index" selection moves to D6, "Generator data" selection moves to B37.

$wb = $excel.ActiveWorkbook

# --- Bus connections: update line data for the two 7500 -> 10000 rows ---
$busConnections = $wb.Worksheets.Item("Bus connections")
$busConnections.Range("E2").Value = 10000
$busConnections.Range("E3").Value = 10000

# --- Generator data: update capacity for the 15000 -> 20000 generator ---
$generatorData = $wb.Worksheets.Item("Generator data")
$generatorData.Range("E4").Value = 20000

# --- Selection bookkeeping (also updates tabSelected / workbook activeTab) ---
$busIndex = $wb.Worksheets.Item("Bus index")
$busIndex.Range("D6").Select()

$generatorData.Range("B37").Select()

# "Bus connections" ends up as the active/selected sheet, with its own selection
# unchanged (still E4), matching the target workbook state.
$busConnections.Activate()
